# Nursery example sheet - add "Owner Name" and "Is Member ? (Yes/No)" columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (inherit the bold header style from the row, same as A1:E1).
$ws.Range("F1").Value = "Owner Name"
$ws.Range("G1").Value = "Is Member ? (Yes/No)"

# Resize the affected columns to (closely) match their content-fit width.
# Column widths below are chosen so the stored column width lands as near as
# possible to Excel's own "best fit" measurement for each header's text.
$ws.Columns.Item(1).ColumnWidth = 5.451822916666667   # A - "Name"
$ws.Columns.Item(2).ColumnWidth = 10.451822916666666  # B - "Contact No."
$ws.Columns.Item(3).ColumnWidth = 7.307291666666667   # C - "Address"
$ws.Columns.Item(6).ColumnWidth = 11.877604166666666  # F - "Owner Name"
$ws.Columns.Item(7).ColumnWidth = 19.877604166666668  # G - "Is Member ? (Yes/No)"

# Move the active selection, as in the committed workbook.
$ws.Range("G11").Select() | Out-Null
